# T460 - joi 7 iulie 2022, 09:52:52 +0300
# Update "Foaie de parcurs" (trip log) for B 151 VGT - iunie 2022 - Alex Bora
# New starting odometer reading and revised daily trip entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: clear every "Locul deplasarii" / "Observatii utilizator" cell so
# none of the old location strings remain referenced anywhere in the sheet.
$ws.Range("C15").Value = ""
$ws.Range("D15").Value = ""
$ws.Range("C16").Value = ""
$ws.Range("D16").Value = ""
$ws.Range("C19").Value = ""
$ws.Range("D19").Value = ""
$ws.Range("C20").Value = ""
$ws.Range("D20").Value = ""
$ws.Range("C21").Value = ""
$ws.Range("D21").Value = ""
$ws.Range("C22").Value = ""
$ws.Range("D22").Value = ""
$ws.Range("C23").Value = ""
$ws.Range("D23").Value = ""
$ws.Range("C27").Value = ""
$ws.Range("D27").Value = ""
$ws.Range("C28").Value = ""
$ws.Range("D28").Value = ""
$ws.Range("C29").Value = ""
$ws.Range("D29").Value = ""
$ws.Range("C30").Value = ""
$ws.Range("D30").Value = ""
$ws.Range("C33").Value = ""
$ws.Range("D33").Value = ""
$ws.Range("C34").Value = ""
$ws.Range("D34").Value = ""
$ws.Range("C35").Value = ""
$ws.Range("D35").Value = ""
$ws.Range("C36").Value = ""
$ws.Range("D36").Value = ""
$ws.Range("C37").Value = ""
$ws.Range("D37").Value = ""
$ws.Range("C40").Value = ""
$ws.Range("D40").Value = ""
$ws.Range("C41").Value = ""
$ws.Range("D41").Value = ""
$ws.Range("C42").Value = ""
$ws.Range("D42").Value = ""
$ws.Range("C43").Value = ""
$ws.Range("D43").Value = ""

# Persist so the now-unreferenced strings actually drop out of the workbook.
$wb.Save()

# --- Step 2: Km initiali (starting odometer) for the period
$ws.Range("B12").Value = 267316

# --- Step 3: re-populate each day's row, in row order, so new/ reused
# location strings come back in the same order the source workbook has them.

# Row 15 (Ziua 2) - unchanged
$ws.Range("B15").Value = 156
$ws.Range("C15").Value = "Cluj-Zalau"
$ws.Range("D15").Value = "Interes Serviciu"

# Row 16 (Ziua 3)
$ws.Range("B16").Value = 47
$ws.Range("C16").Value = "Cluj-Cluj"
$ws.Range("D16").Value = "Interes Serviciu"

# Row 19 (Ziua 6)
$ws.Range("B19").Value = 30
$ws.Range("C19").Value = "Acasa-Birou"
$ws.Range("D19").Value = " "

# Row 20 (Ziua 7)
$ws.Range("B20").Value = 121
$ws.Range("C20").Value = "Cluj-Turda"
$ws.Range("D20").Value = "Interes Serviciu"

# Row 21 (Ziua 8)
$ws.Range("B21").Value = 101
$ws.Range("C21").Value = "Cluj-Dej"
$ws.Range("D21").Value = "Interes Serviciu"

# Row 22 (Ziua 9) - unchanged
$ws.Range("B22").Value = 156
$ws.Range("C22").Value = "Cluj-Zalau"
$ws.Range("D22").Value = "Interes Serviciu"

# Row 23 (Ziua 10)
$ws.Range("B23").Value = 92
$ws.Range("C23").Value = "Cluj-Bontida"
$ws.Range("D23").Value = "Interes Serviciu"

# Row 27 (Ziua 14)
$ws.Range("B27").Value = 257
$ws.Range("C27").Value = "Cluj-Bistrita"
$ws.Range("D27").Value = "Interes Serviciu"

# Row 28 (Ziua 15)
$ws.Range("B28").Value = 101
$ws.Range("C28").Value = "Cluj-Dej"
$ws.Range("D28").Value = "Interes Serviciu"

# Row 29 (Ziua 16)
$ws.Range("B29").Value = 85
$ws.Range("C29").Value = "Cluj-Apahida"
$ws.Range("D29").Value = "Interes Serviciu"

# Row 30 (Ziua 17)
$ws.Range("B30").Value = 30
$ws.Range("C30").Value = "Acasa-Birou"
$ws.Range("D30").Value = " "

# Row 33 (Ziua 20)
$ws.Range("B33").Value = 30
$ws.Range("C33").Value = "Acasa-Birou"
$ws.Range("D33").Value = " "

# Row 34 (Ziua 21)
$ws.Range("B34").Value = 152
$ws.Range("C34").Value = "Cluj-Cmp. Turzii"
$ws.Range("D34").Value = "Interes Serviciu"

# Row 35 (Ziua 22)
$ws.Range("B35").Value = 257
$ws.Range("C35").Value = "Cluj-Bistrita"
$ws.Range("D35").Value = "Interes Serviciu"

# Row 36 (Ziua 23)
$ws.Range("B36").Value = 30
$ws.Range("C36").Value = "Acasa-Birou"
$ws.Range("D36").Value = " "

# Row 37 (Ziua 24)
$ws.Range("B37").Value = 257
$ws.Range("C37").Value = "Cluj-Bistrita"
$ws.Range("D37").Value = "Interes Serviciu"

# Row 40 (Ziua 27)
$ws.Range("B40").Value = 421
$ws.Range("C40").Value = "Cluj-Satu-Mare"
$ws.Range("D40").Value = "Interes Serviciu"

# Row 41 (Ziua 28)
$ws.Range("B41").Value = 30
$ws.Range("C41").Value = "Acasa-Birou"
$ws.Range("D41").Value = " "

# Row 42 (Ziua 29)
$ws.Range("B42").Value = 30
$ws.Range("C42").Value = "Acasa-Birou"
$ws.Range("D42").Value = " "

# Row 43 (Ziua 30)
$ws.Range("B43").Value = 121
$ws.Range("C43").Value = "Cluj-Turda"
$ws.Range("D43").Value = "Interes Serviciu"

# --- Step 4: totals
$ws.Range("B44").Value = 2504
$ws.Range("B45").Value = 269820
